{"js": "// The edit collapses several multi-run paragraphs into single runs (no\n// visible text/formatting change), makes \"Figure 5:\" bold (matching the\n// other figure headers), expands \"Files : \" with the actual file names,\n// and inserts four new informational paragraphs describing the new\n// \"Figure 5\" dataset (Model / Genes / quencher / sample-group column),\n// including a bookmark around the gene list, as in the other model\n// description blocks already present in the document.\n\nasync function mergeParagraphRuns(searchText) {\n  const results = context.document.body.search(searchText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    return;\n  }\n  const range = results.items[0];\n  // Re-writing the full paragraph text via Replace collapses every run\n  // in the range into a single run while preserving the formatting of\n  // the first run - exactly what the diff shows for these paragraphs.\n  range.insertText(searchText, \"Replace\");\n  await context.sync();\n}\n\n// 1) Collapse the three \"Control ... order\" caption paragraphs (each was\n//    split across 3 runs: \"prefix [\", \"option text\", \"]\") into one run.\nawait mergeParagraphRuns(\"Control cell lines and stages order ------ [used in paper]\");\nawait mergeParagraphRuns(\"Control order of stages -----  [another option]\");\nawait mergeParagraphRuns(\"Control cell lines order ----- [also works]\");\n\n// 2) Collapse the \"Statistic relative\" selection paragraphs (each split\n//    into 2-3 runs) into a single run.\nawait mergeParagraphRuns(\"select: Group names are within the sample name \u2013 group names are: D0,D7\");\nawait mergeParagraphRuns(\"select for repeated measures: No --- measures are independent \");\nawait mergeParagraphRuns(\"select for normal distribution: Yes --- Parametric tests\");\n\n// 3) Make the \"Figure 5:\" heading bold, matching the other figure\n//    headings (Figure 2/3/4) in the document.\n{\n  const results = context.document.body.search(\"Figure 5:\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  const range = results.items[0];\n  const para = range.paragraphs.getFirst();\n  para.font.bold = true;\n  await context.sync();\n}\n\n// 4) Update \"Files : \" with the actual file names used for this figure.\n{\n  const results = context.document.body.search(\"Files : \", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  const range = results.items[0];\n  range.insertText(\"Files : B2M.csv, NRXN3.csv\", \"Replace\");\n  await context.sync();\n}\n\n// 5) Insert the new paragraphs describing this dataset right after the\n//    \"Files : ...\" paragraph and before \"Targets for normalization...\".\n{\n  const results = context.document.body.search(\"Files : B2M.csv, NRXN3.csv\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  const filesRange = results.items[0];\n  let anchorPara = filesRange.paragraphs.getFirst();\n\n  anchorPara = anchorPara.insertParagraph(\"Model: Absolute\", \"After\");\n  await context.sync();\n\n  const genesPara = anchorPara.insertParagraph(\"Genes if file name only: B2M,NRXN3\", \"After\");\n  await context.sync();\n\n  // Bookmark the gene list (\"B2M,NRXN3\") within the newly inserted paragraph.\n  const geneResults = genesPara.search(\"B2M,NRXN3\", { matchCase: true });\n  geneResults.load(\"items\");\n  await context.sync();\n  geneResults.items[0].insertBookmark(\"__DdeLink__9670_3289051397\");\n  await context.sync();\n\n  const quencherPara = genesPara.insertParagraph(\"Name of quencher:TMR\", \"After\");\n  await context.sync();\n\n  quencherPara.insertParagraph(\n    \"Name of the sample group in the task or content column: sample\",\n    \"After\"\n  );\n  await context.sync();\n}\n\n// 6) Collapse the \"Sampel order (cell lines, time points) \u2013 in this case\n//    brain regions and mice\" paragraph (2 runs) into a single run.\nawait mergeParagraphRuns(\"Sampel order (cell lines, time points) \u2013 in this case brain regions and mice\");\n", "ps1": "# The edit collapses several multi-run paragraphs into single runs (no\n# visible text/formatting change), makes \"Figure 5:\" bold (matching the\n# other figure headers), expands \"Files : \" with the actual file names,\n# and inserts four new informational paragraphs describing the new\n# \"Figure 5\" dataset (Model / Genes / quencher / sample-group column),\n# including a bookmark around the gene list, as in the other model\n# description blocks already present in the document.\n\nfunction Get-ParagraphByText($doc, $targetText) {\n    $paras = $doc.Paragraphs\n    for ($i = 1; $i -le $paras.Count; $i++) {\n        $p = $paras.Item($i)\n        $t = $p.Range.Text\n        $t = $t.TrimEnd([char]13, [char]7)\n        if ($t -eq $targetText) {\n            return $p\n        }\n    }\n    return $null\n}\n\nfunction Merge-ParagraphRuns($doc, $targetText) {\n    $p = Get-ParagraphByText $doc $targetText\n    if ($p -ne $null) {\n        # Searching+replacing the paragraph's own range with its current\n        # text collapses every run inside it into a single run (keeping\n        # the first run's formatting) - this is exactly what the diff\n        # shows for these paragraphs.\n        $find = $p.Range.Find\n        $find.Execute($targetText, $false, $false, $false, $false, $false, $true, 1, $false, $targetText, 2) | Out-Null\n    }\n}\n\n$d = $word.ActiveDocument\n\n# 1) Collapse the three \"Control ... order\" caption paragraphs (each was\n#    split across 3 runs: \"prefix [\", \"option text\", \"]\") into one run.\nMerge-ParagraphRuns $d \"Control cell lines and stages order ------ [used in paper]\"\nMerge-ParagraphRuns $d \"Control order of stages -----  [another option]\"\nMerge-ParagraphRuns $d \"Control cell lines order ----- [also works]\"\n\n# 2) Collapse the \"Statistic relative\" selection paragraphs (each split\n#    into 2-3 runs) into a single run.\nMerge-ParagraphRuns $d \"select: Group names are within the sample name \u2013 group names are: D0,D7\"\nMerge-ParagraphRuns $d \"select for repeated measures: No --- measures are independent \"\nMerge-ParagraphRuns $d \"select for normal distribution: Yes --- Parametric tests\"\n\n# 3) Make the \"Figure 5:\" heading bold, matching the other figure\n#    headings (Figure 2/3/4) in the document.\n$figure5 = Get-ParagraphByText $d \"Figure 5:\"\n$figure5.Range.Font.Bold = $true\n$figure5.Range.Font.BoldBi = $true\n\n# 4) Update \"Files : \" with the actual file names used for this figure.\n$filesPara = Get-ParagraphByText $d \"Files : \"\n$find = $filesPara.Range.Find\n$find.Execute(\"Files : \", $false, $false, $false, $false, $false, $true, 1, $false, \"Files : B2M.csv, NRXN3.csv\", 2) | Out-Null\n\n# 5) Insert the new paragraphs describing this dataset right after the\n#    \"Files : ...\" paragraph and before \"Targets for normalization...\".\n$filesPara = Get-ParagraphByText $d \"Files : B2M.csv, NRXN3.csv\"\n$idx = $filesPara.Index\n\n$filesPara.Range.InsertParagraphAfter()\n$modelPara = $d.Paragraphs.Item($idx + 1)\n$modelPara.Range.Text = \"Model: Absolute\"\n\n$modelPara.Range.InsertParagraphAfter()\n$genesPara = $d.Paragraphs.Item($idx + 2)\n$genesPara.Range.Text = \"Genes if file name only: B2M,NRXN3\"\n\n# Bookmark the gene list (\"B2M,NRXN3\") within the newly inserted paragraph.\n$genesRange = $genesPara.Range\n$bmFind = $genesRange.Find\n$bmFind.Execute(\"B2M,NRXN3\") | Out-Null\n$d.Bookmarks.Add(\"__DdeLink__9670_3289051397\", $genesRange)\n\n$genesPara.Range.InsertParagraphAfter()\n$quencherPara = $d.Paragraphs.Item($idx + 3)\n$quencherPara.Range.Text = \"Name of quencher:TMR\"\n\n$quencherPara.Range.InsertParagraphAfter()\n$samplePara = $d.Paragraphs.Item($idx + 4)\n$samplePara.Range.Text = \"Name of the sample group in the task or content column: sample\"\n\n# 6) Collapse the \"Sampel order (cell lines, time points) \u2013 in this case\n#    brain regions and mice\" paragraph (2 runs) into a single run.\nMerge-ParagraphRuns $d \"Sampel order (cell lines, time points) \u2013 in this case brain regions and mice\"\n\nWrite-Output \"Done\"\n"}
